# Updates the cryptos list (price + 1h volume change columns, plus a
# reordering of the PEPE / WrappedeETH rows) to match the refreshed data
# pulled by the "Updated cryptos list ... with GitHub Actions" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "62.392.19"
$ws.Range("E2").Value = "  -1.42%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.459.41"
$ws.Range("E3").Value = "  +0.23%  "

# --- Row 4: TetherUSD ---
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'582.72"
$ws.Range("E5").Value = "  +1.66%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "'144.06"
$ws.Range("E6").Value = "  -1.66%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.02%  "

# --- Row 8: XRP ---
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  -1.07%  "

# --- Row 9: LidoStakedEther ---
$ws.Range("D9").Value = "2.456.75"
$ws.Range("E9").Value = "  +0.22%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -3.48%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  +2.19%  "

# --- Row 12: Toncoin ---
$ws.Range("E12").Value = "  -0.87%  "

# --- Row 13: Cardano ---
$ws.Range("E13").Value = "  -2.66%  "

# --- Row 14: Avalanche ---
$ws.Range("D14").Value = "'26.65"
$ws.Range("E14").Value = "  -1.52%  "

# --- Row 15: ShibaInu ---
$ws.Range("E15").Value = "  -3.75%  "

# --- Row 16: WrappedliquidstakedEther2.0 ---
$ws.Range("D16").Value = "2.824.10"
$ws.Range("E16").Value = "  -2.67%  "

# --- Row 17: WrappedBTC ---
$ws.Range("D17").Value = "62.164.91"
$ws.Range("E17").Value = "  -1.78%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "2.450.31"
$ws.Range("E18").Value = "  +0.26%  "

# --- Row 19: Chainlink ---
$ws.Range("D19").Value = "'10.94"
$ws.Range("E19").Value = "  -3.48%  "

# --- Row 20 ---
$ws.Range("D20").Value = "'7.20"
$ws.Range("E20").Value = "  -1.46%  "

# --- Row 21 ---
$ws.Range("D21").Value = "'331.12"
$ws.Range("E21").Value = "  +0.56%  "

# --- Row 22 ---
$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  -2.28%  "

# --- Row 23 ---
$ws.Range("D23").Value = "'2.01"
$ws.Range("E23").Value = "  -3.31%  "

# --- Row 24 ---
$ws.Range("E24").Value = "  +0.09%  "

# --- Row 25: Litecoin ---
$ws.Range("D25").Value = "'66.08"
$ws.Range("E25").Value = "  +0.75%  "

# --- Row 26: Aptos ---
$ws.Range("E26").Value = "  +6.39%  "

# --- Row 27: Bittensor ---
$ws.Range("D27").Value = "'625.78"
$ws.Range("E27").Value = "  +1.59%  "

# --- Row 28 & 29: PEPE and WrappedeETH swap places ---
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0961"
$ws.Range("E28").Value = "  -6.32%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.542.32"
$ws.Range("E29").Value = "  -0.87%  "

# --- Row 30: Binance-PegBSC-USD ---
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.45%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -4.04%  "

# --- Row 32 ---
$ws.Range("E32").Value = "  -2.11%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  +0.81%  "

# --- Row 34 ---
$ws.Range("E34").Value = "  -0.67%  "

# --- Row 35: NEARProtocol ---
$ws.Range("D35").Value = "'4.94"
$ws.Range("E35").Value = "  -4.95%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  +0.17%  "

# --- Row 37: ImmutableX ---
$ws.Range("D37").Value = "'1.43"
$ws.Range("E37").Value = "  -6.33%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  -0.15%  "

# --- Row 39 ---
$ws.Range("E39").Value = "  -0.92%  "

# --- Row 40: Monero ---
$ws.Range("D40").Value = "'149.72"
$ws.Range("E40").Value = "  +2.00%  "

# --- Row 41: EthereumClassic ---
$ws.Range("D41").Value = "'18.42"
$ws.Range("E41").Value = "  -2.34%  "

# --- Row 42 ---
$ws.Range("E42").Value = "  -2.15%  "

# --- Row 43: OKB ---
$ws.Range("D43").Value = "'42.52"
$ws.Range("E43").Value = "  +1.73%  "

# --- Row 45: dogwifhat ---
$ws.Range("D45").Value = "'2.48"
$ws.Range("E45").Value = "  -4.67%  "

# --- Row 46: Aave ---
$ws.Range("D46").Value = "'144.01"
$ws.Range("E46").Value = "  -3.13%  "

# --- Row 47 ---
$ws.Range("E47").Value = "  -3.04%  "

# --- Row 48: Hedera ---
$ws.Range("D48").Value = "'0.0528"
$ws.Range("E48").Value = "  -1.04%  "

# --- Row 49 ---
$ws.Range("E49").Value = "  +0.43%  "

# --- Row 50: InjectiveProtocol ---
$ws.Range("D50").Value = "'19.73"
$ws.Range("E50").Value = "  -6.85%  "

# --- Row 51 ---
$ws.Range("E51").Value = "  +9.18%  "
